# Populate the previously-empty "Clothing" column (G) for rows 2-90
# with comma-separated clothing-label predictions from the trained model.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$clothingValues = @(
    'Kaftan,Tee',
    'Kaftan,Jumpsuit',
    'Jumpsuit,Top',
    'Jumpsuit,Blouse',
    'Trunks,Caftan',
    'Halter,Jumpsuit',
    'Jumpsuit,Dress',
    'Halter,Jumpsuit',
    'Jumpsuit,Blouse',
    'Trunks,Jumpsuit',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Dress',
    'Halter,Trunks',
    'Blouse,Jumpsuit',
    'Blazer,Jumpsuit',
    'Kaftan,Jumpsuit',
    'Kaftan,Jumpsuit',
    'Blouse,Kaftan',
    'Jumpsuit,Blouse',
    'Jumpsuit,Halter',
    'Blouse,Blazer',
    'Blouse,Kaftan',
    'Blouse,Kaftan',
    'Halter,Blazer',
    'Jumpsuit,Jodhpurs',
    'Jumpsuit,Dress',
    'Sweatpants,Blazer',
    'Jumpsuit,Kaftan',
    'Blouse,Caftan',
    'Jumpsuit,Blouse',
    'Blouse,Jumpsuit',
    'Jumpsuit,Halter',
    'Jumpsuit,Parka',
    'Trunks,Jodhpurs',
    'Jumpsuit,Halter',
    'Jumpsuit,Trunks',
    'Blouse,Jumpsuit',
    'Jumpsuit,Blouse',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Blouse',
    'Jumpsuit,Blouse',
    'Blouse,Caftan',
    'Blouse,Caftan',
    'Parka,Gauchos',
    'Blouse,Trunks',
    'Blouse,Tee',
    'Jumpsuit,Dress',
    'Jumpsuit,Halter',
    'Jumpsuit,Halter',
    'Parka,Blouse',
    'Blazer,Jumpsuit',
    'Jumpsuit,Blouse',
    'Jodhpurs,Jumpsuit',
    'Blouse,Parka',
    'Blouse,Jumpsuit',
    'Caftan,Jumpsuit',
    'Jumpsuit,Blouse',
    'Jumpsuit,Dress',
    'Jodhpurs,Trunks',
    'Jumpsuit,Kaftan',
    'Blouse,Jumpsuit',
    'Caftan,Jumpsuit',
    'Jumpsuit,Kaftan',
    'Cutoffs,Coat',
    'Jumpsuit,Kaftan',
    'Tee,Kaftan',
    'Jumpsuit,Parka',
    'Halter,Blazer',
    'Blouse,Jumpsuit',
    'Jumpsuit,Chinos',
    'Halter,Blazer',
    'Jumpsuit,Parka',
    'Jumpsuit,Parka',
    'Parka,Jumpsuit',
    'Caftan,Trunks',
    'Jumpsuit,Kaftan',
    'Dress,Parka',
    'Halter,Jumpsuit',
    'Tee,Jumpsuit',
    'Jumpsuit,Tee',
    'Trunks,Jodhpurs',
    'Halter,Jumpsuit',
    'Caftan,Halter',
    'Tee,Jumpsuit',
    'Parka,Sweatpants',
    'Trunks,Jumpsuit',
    'Jumpsuit,Tee',
    'Parka,Blouse',
    'Jumpsuit,Parka'
)

for ($i = 0; $i -lt $clothingValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $clothingValues[$i]
}
